$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D5").Value = 0.26
Write-Host "done D5"
